$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '64.660.19'
$ws.Range('E2').Value = '  -1.72%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '3.419.74'
$ws.Range('E3').Value = '  -1.90%  '
$ws.Range('E4').Value = '  +0.01%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '572.30'
$ws.Range('E5').Value = '  -1.44%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '157.81'
$ws.Range('E6').Value = '  -1.95%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.611'
$ws.Range('E7').Value = '  +0.85%  '
$ws.Range('E8').Value = '  +0.04%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '3.420.22'
$ws.Range('E9').Value = '  -1.90%  '
$ws.Range('E10').Value = '  -1.71%  '
$ws.Range('E11').Value = '  -2.46%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.439'
$ws.Range('E12').Value = '  -1.88%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '4.013.94'
$ws.Range('E13').Value = '  -1.75%  '
$ws.Range('E14').Value = '  -0.29%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.0000187'
$ws.Range('E15').Value = '  -4.43%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '27.65'
$ws.Range('E16').Value = '  -4.06%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '64.710.84'
$ws.Range('E17').Value = '  -1.54%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '3.453.90'
$ws.Range('E18').Value = '  -0.74%  '
$ws.Range('E19').Value = '  -2.24%  '
$ws.Range('E20').Value = '  -3.70%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '379.84'
$ws.Range('E21').Value = '  -2.80%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '7.98'
$ws.Range('E22').Value = '  -3.26%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '0.546'
$ws.Range('E23').Value = '  -1.18%  '
$ws.Range('E24').Value = '  -0.17%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '71.99'
$ws.Range('E25').Value = '  -2.05%  '
$ws.Range('E26').Value = '  -5.68%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '9.97'
$ws.Range('E27').Value = '  +1.80%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '0.177'
$ws.Range('E28').Value = '  -0.59%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '1.01'
$ws.Range('E29').Value = '  +0.79%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '1.46'
$ws.Range('E30').Value = '  +1.17%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '6.15'
$ws.Range('E31').Value = '  -3.84%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '2.00'
$ws.Range('E32').Value = '  -2.86%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '23.19'
$ws.Range('E33').Value = '  -2.47%  '
$ws.Range('E34').Value = '  -1.17%  '
$ws.Range('E35').Value = '  +1.83%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '160.31'
$ws.Range('E36').Value = '  -1.66%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '1.90'
$ws.Range('E37').Value = '  -3.09%  '
$ws.Range('B38').Value = 'Hedera'
$ws.Range('C38').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.0754'
$ws.Range('E38').Value = '  -2.65%  '
$ws.Range('B39').Value = 'Maker'
$ws.Range('C39').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '2.895.25'
$ws.Range('E39').Value = '  -6.22%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '6.68'
$ws.Range('E40').Value = '  +2.77%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '26.32'
$ws.Range('E41').Value = '  -3.29%  '
$ws.Range('E42').Value = '  +0.54%  '
$ws.Range('E43').Value = '  -0.26%  '
$ws.Range('E44').Value = '  -2.31%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.770'
$ws.Range('E45').Value = '  -1.29%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '25.73'
$ws.Range('E46').Value = '  -0.31%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '316.65'
$ws.Range('E47').Value = '  +1.82%  '
$ws.Range('E48').Value = '  -0.99%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '1.06'
$ws.Range('E49').Value = '  -5.04%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.107'
$ws.Range('E50').Value = '  -1.81%  '
$ws.Range('E51').Value = '  -3.12%  '
